$wb = $excel.ActiveWorkbook

# --- adc2ver sheet: update the resistor-divider / ADC version table ---
$ws = $wb.Worksheets.Item("adc2ver")

# Update the input resistor values (column C) for rows 4-6, and row 7 value
# (row 7's value also loses its explicit number-format style, matching the
# original edit where the cell format was cleared before re-entering data).
$ws.Range("C4").Value = 4700
$ws.Range("C5").Value = 10000
$ws.Range("C6").Value = 47000
$ws.Range("C7").ClearFormats() | Out-Null
$ws.Range("C7").Value = 100000

# Re-enter the dependent formulas so they recompute off the new inputs.
$ws.Range("D2").Formula = "=1.8*C2/(C2+B2)"
$ws.Range("D3").Formula = "=1.8*C3/(C3+B3)"
$ws.Range("D4").Formula = "=1.8*C4/(C4+B4)"
$ws.Range("D5").Formula = "=1.8*C5/(C5+B5)"
$ws.Range("D6").Formula = "=1.8*C6/(C6+B6)"
$ws.Range("D7").Formula = "=1.8*C7/(C7+B7)"

$ws.Range("F2").Formula = "=INT(E2/1.8*2^12)"
$ws.Range("F3").Formula = "=INT(E3/1.8*2^12)"
$ws.Range("F4").Formula = "=INT(E4/1.8*2^12)"
$ws.Range("F5").Formula = "=INT(E5/1.8*2^12)"
$ws.Range("F6").Formula = "=INT(E6/1.8*2^12)"
$ws.Range("F7").Formula = "=INT(E7/1.8*2^12)"

# Widen the new/adjusted columns (C, D, E) to fit the larger numbers.
$ws.Columns.Item(3).ColumnWidth = 6.466666666666667
$ws.Range("D1:E1").ColumnWidth = 11.666666666666666

# Make adc2ver the active sheet/tab, with D5 as the selected cell - this
# also clears tabSelected on whichever sheet was previously active (pin).
$ws.Activate() | Out-Null
$ws.Range("D5").Select() | Out-Null

Write-Host "done"
